{"js": "// Replace the three-digit-by-one-digit multiplication answers in the\n// practice-sheet table with the new set of problems/answers.\n// Each (old, new) pair below is a unique, whole-cell text value, so a\n// plain text search-and-replace (matching the whole cell text) is safe.\nconst replacements = [\n  [\"378\u00d79=3402\", \"238\u00d72=476\"],\n  [\"633\u00d75=3165\", \"250\u00d72=500\"],\n  [\"787\u00d79=7083\", \"680\u00d73=2040\"],\n  [\"207\u00d78=1656\", \"731\u00d72=1462\"],\n  [\"805\u00d73=2415\", \"894\u00d74=3576\"],\n  [\"495\u00d79=4455\", \"253\u00d75=1265\"],\n  [\"922\u00d78=7376\", \"236\u00d78=1888\"],\n  [\"136\u00d78=1088\", \"996\u00d75=4980\"],\n  [\"982\u00d79=8838\", \"966\u00d79=8694\"],\n  [\"206\u00d78=1648\", \"269\u00d77=1883\"],\n  [\"459\u00d77=3213\", \"526\u00d78=4208\"],\n  [\"271\u00d78=2168\", \"356\u00d76=2136\"],\n  [\"361\u00d73=1083\", \"297\u00d79=2673\"],\n  [\"752\u00d72=1504\", \"395\u00d77=2765\"],\n  [\"457\u00d74=1828\", \"687\u00d77=4809\"],\n  [\"871\u00d79=7839\", \"467\u00d76=2802\"],\n  [\"310\u00d78=2480\", \"923\u00d76=5538\"],\n  [\"762\u00d72=1524\", \"202\u00d75=1010\"],\n  [\"117\u00d79=1053\", \"792\u00d74=3168\"],\n  [\"455\u00d79=4095\", \"921\u00d72=1842\"],\n  [\"445\u00d75=2225\", \"426\u00d77=2982\"],\n  [\"706\u00d78=5648\", \"141\u00d76=846\"],\n  [\"905\u00d79=8145\", \"337\u00d78=2696\"],\n  [\"491\u00d74=1964\", \"261\u00d77=1827\"],\n  [\"356\u00d74=1424\", \"768\u00d74=3072\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the three-digit-by-one-digit multiplication answers in the\n# practice-sheet table with the new set of problems/answers.\n# Each (old, new) pair below is a unique, whole-cell text value, so a\n# plain Find/Replace on the exact text is safe (single hit each).\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"378\u00d79=3402\", \"238\u00d72=476\"),\n    @(\"633\u00d75=3165\", \"250\u00d72=500\"),\n    @(\"787\u00d79=7083\", \"680\u00d73=2040\"),\n    @(\"207\u00d78=1656\", \"731\u00d72=1462\"),\n    @(\"805\u00d73=2415\", \"894\u00d74=3576\"),\n    @(\"495\u00d79=4455\", \"253\u00d75=1265\"),\n    @(\"922\u00d78=7376\", \"236\u00d78=1888\"),\n    @(\"136\u00d78=1088\", \"996\u00d75=4980\"),\n    @(\"982\u00d79=8838\", \"966\u00d79=8694\"),\n    @(\"206\u00d78=1648\", \"269\u00d77=1883\"),\n    @(\"459\u00d77=3213\", \"526\u00d78=4208\"),\n    @(\"271\u00d78=2168\", \"356\u00d76=2136\"),\n    @(\"361\u00d73=1083\", \"297\u00d79=2673\"),\n    @(\"752\u00d72=1504\", \"395\u00d77=2765\"),\n    @(\"457\u00d74=1828\", \"687\u00d77=4809\"),\n    @(\"871\u00d79=7839\", \"467\u00d76=2802\"),\n    @(\"310\u00d78=2480\", \"923\u00d76=5538\"),\n    @(\"762\u00d72=1524\", \"202\u00d75=1010\"),\n    @(\"117\u00d79=1053\", \"792\u00d74=3168\"),\n    @(\"455\u00d79=4095\", \"921\u00d72=1842\"),\n    @(\"445\u00d75=2225\", \"426\u00d77=2982\"),\n    @(\"706\u00d78=5648\", \"141\u00d76=846\"),\n    @(\"905\u00d79=8145\", \"337\u00d78=2696\"),\n    @(\"491\u00d74=1964\", \"261\u00d77=1827\"),\n    @(\"356\u00d74=1424\", \"768\u00d74=3072\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for '$oldText'\"\n    }\n}\n"}
